# Venezuela Primera Division - base update (13-06-2024 19:35)
# The underlying match records for several fixtures were re-keyed (ids
# reshuffled/corrected) which, in the exported row-per-match sheet, shows
# up as the odds/result data for a handful of rows being rotated among
# the rows that share the same match date. Column A (row/display index)
# and columns C/D (Div/Date) are untouched; only B and E:AD move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Rows($Rows) {
    # Captures B:AD for every row in $Rows, then writes row[i]'s new data
    # from row[i+1]'s old data (cyclically): new(Rows[i]) = old(Rows[i+1]).
    $snapshots = @{}
    foreach ($r in $Rows) {
        $rng = "B$r`:AD$r"
        $snapshots[$r] = $ws.Range($rng).Value2
    }
    $count = $Rows.Length
    for ($i = 0; $i -lt $count; $i++) {
        $target = $Rows[$i]
        $source = $Rows[($i + 1) % $count]
        $rng = "B$target`:AD$target"
        $ws.Range($rng).Value = $snapshots[$source]
    }
}

# Rows 94, 96, 97 (ids 6236251/6236252/6236254) -- 3-way rotation
Rotate-Rows @(94, 96, 97)

# Rows 102, 103 (ids 6236616/6236615) -- simple swap
Rotate-Rows @(102, 103)

# Rows 114, 115 (ids 7352251/7352250) -- simple swap
Rotate-Rows @(114, 115)

# Rows 162, 163 (ids 7952905/7952893) -- simple swap
Rotate-Rows @(162, 163)

# Rows 204, 205, 206 (ids 7977880/7977876/7977386) -- 3-way rotation
Rotate-Rows @(204, 205, 206)
